# Edit: Fix "Basic transfers" capitalization and insert new "Burst Transfers"
# test-plan rows (with sub-cases) into the Stimulus sheet, pushing the final
# "Cross Feature" row down. A light top border is used on the first inserted
# row's Tests cell to visually separate burst sub-cases from the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stimulus")

# --- Insert 10 new rows before current row 11 (the "Cross Feature" row) ---
# After insertion, old row 11 becomes row 21, and the 10 new rows occupy 11..20.
$insertRange = $ws.Range("B11:D20")
$insertRange.Insert(-4162) # xlShiftDown

# Style references used below, matching the existing table formatting:
#  - column B (Test No) numeric cells use the centered/bordered style already
#    applied by the sheet (copied down automatically with Insert above row 11..20
#    inherited formatting from row 11's original style since Insert copies the
#    format of the row above by default in Excel; we explicitly reapply below
#    to be safe).

# Reapply formatting for B, C, D across rows 11-21 to be certain it matches
# the rest of the table (same look as rows 4-10).
$ws.Range("B11:B21").Font.Bold = $false
$ws.Range($ws.Cells.Item(4,2), $ws.Cells.Item(4,2)).Copy() | Out-Null
$ws.Range("B11:B21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range($ws.Cells.Item(4,3), $ws.Cells.Item(4,3)).Copy() | Out-Null
$ws.Range("C11:C21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range($ws.Cells.Item(4,4), $ws.Cells.Item(4,4)).Copy() | Out-Null
$ws.Range("D11:D21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Fill in the Test No (column B) values 8..18 ---
for ($i = 0; $i -le 10; $i++) {
    $ws.Cells.Item(11 + $i, 2).Value = 8 + $i
}

# --- Column C (Coverage Plan) ---
# Rows 11-20: "Burst Transfers"; row 21 keeps "Cross Feature"
for ($r = 11; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "Burst Transfers"
}
$ws.Cells.Item(21, 3).Value = "Cross Feature"

# --- Column D (Tests) ---
$ws.Cells.Item(11, 4).Value = "Incrementing Burst (INCR): Check that addresses increment correctly and stay within the subordinate" + [char]0x2019 + "s addressable range."
$ws.Cells.Item(12, 4).Value = "WRAP4 (HBURST = 0b010): Verify that the address wraps correctly at the boundary"
$ws.Cells.Item(13, 4).Value = "WRAP8 (HBURST = 0b100)"
$ws.Cells.Item(14, 4).Value = "WRAP16 (HBURST = 0b110)"
$ws.Cells.Item(15, 4).Value = "INCR4 (HBURST = 0b011)"
$ws.Cells.Item(16, 4).Value = "INCR8 (HBURST = 0b101)"
$ws.Cells.Item(17, 4).Value = "INCR16 (HBURST = 0b111)"
$ws.Cells.Item(18, 4).Value = "Edge Cases:Burst with minimum transfer size (HSIZE = BYTE)"
$ws.Cells.Item(19, 4).Value = "Edge Cases:Burst with maximum transfer size (HSIZE = WORD)"
$ws.Cells.Item(20, 4).Value = "Edge Cases:Early burst termination"
$ws.Cells.Item(21, 4).Value = "Multiple Write + Multiple Read txns with hsize, haddr same for each set of txns, with other fields randomize."

# Fix capitalization on existing row 4 ("Basic transfers" -> "Basic Transfers").
# Done after the "Burst Transfers" strings above so the shared-string table
# ends up in the same append order as the authored workbook.
$ws.Cells.Item(4, 3).Value = "Basic Transfers"

# --- Drop the bottom border on D12's Tests cell (it already has thin left/ ---
# --- right/top borders from the paste above), separating the first burst  ---
# --- sub-case row from the "Incrementing Burst" summary row above it.     ---
# NOTE: use the plain numeric Borders collection index (9 = bottom edge)
# rather than the named xlEdgeBottom constant - indexing with the negative
# named constant here clears sibling edges too.
$ws.Cells.Item(12, 4).Borders.Item(9).LineStyle = -4142 # xlLineStyleNone

# --- Update dimension / selection bookkeeping to match final sheet extent ---
$ws.Range("D30").Select() | Out-Null
